$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.504.93'
$ws.Range('E2').Value = '  -2.81%  '

$ws.Range('D3').Value = '3.431.63'
$ws.Range('E3').Value = '  -1.25%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '581.09'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.71%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '171.65'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.63%  '

$ws.Range('E7').Value = '  +0.03%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.595'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.31%  '

$ws.Range('D9').Value = '3.429.05'
$ws.Range('E9').Value = '  -1.31%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.129'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -5.69%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '6.89'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.32%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.407'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -4.29%  '

$ws.Range('D13').Value = '4.027.90'
$ws.Range('E13').Value = '  -1.32%  '

$ws.Range('E14').Value = '  +0.83%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '28.83'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -9.89%  '

$ws.Range('D16').Value = '65.637.84'
$ws.Range('E16').Value = '  -2.57%  '

$ws.Range('E17').Value = '  -3.86%  '

$ws.Range('D18').Value = '3.433.59'
$ws.Range('E18').Value = '  -1.30%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '5.90'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -3.70%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.80'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.56%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '364.82'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -6.25%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '7.63'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -3.63%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '72.62'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.55%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  -0.99%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.0000121'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.59%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.74'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -4.07%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.177'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.75%  '

$ws.Range('E29').Value = '  +0.05%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '23.63'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -3.86%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.97'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -3.20%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '5.67'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -6.21%  '

$ws.Range('E33').Value = '  +0.02%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.29'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -6.96%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '7.00'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -3.64%  '

$ws.Range('E36').Value = '  -3.28%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '160.61'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.23%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '28.85'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.14%  '

$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.879'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.62%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.58'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.07%  '

$ws.Range('D41').Value = '2.756.87'
$ws.Range('E41').Value = '  +1.19%  '

$ws.Range('E42').Value = '  -5.46%  '

$ws.Range('E43').Value = '  -2.57%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.39'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -3.38%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0676'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.63%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '39.86'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -3.65%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '24.01'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -7.38%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0287'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -3.37%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '321.45'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -3.25%  '

$ws.Range('E50').Value = '  -3.47%  '

$ws.Range('E51').Value = '  -1.41%  '
